$d = $word.ActiveDocument

# Locate the paragraph that contains the sentence being edited.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*those drive high sales*") {
        $targetPara = $p
        break
    }
}

$pStart = $targetPara.Range.Start
$pText = $targetPara.Range.Text

$oldPhrase = "those drive"
$newPhrase = "contributed"

$phraseStart = $pStart + $pText.IndexOf($oldPhrase)
$phraseEnd = $phraseStart + $oldPhrase.Length

# Boundary between the original run's text ("...campaigns. R") and the
# following run ("esponsible...") so we can re-split it out after the
# engine coalesces same-formatted runs on edit.
$sentenceMarker = "campaigns. R"
$afterRIdx = $pStart + $pText.IndexOf($sentenceMarker) + $sentenceMarker.Length

# 1) Swap the phrase text (same length, so offsets after it don't move).
$phraseRange = $d.Range($phraseStart, $phraseEnd)
$phraseRange.Text = $newPhrase

$afterPhrase = $phraseStart + $newPhrase.Length

# 2) Re-establish a run boundary right after "contributed" by touching the
#    formatting (set then restore the point size) on just that word.
$origSize = 8
$wordRange = $d.Range($phraseStart, $afterPhrase)
$wordRange.Font.Size = $origSize + 1
$wordRange2 = $d.Range($phraseStart, $afterPhrase)
$wordRange2.Font.Size = $origSize

# 3) Re-establish a run boundary right after "...campaigns. R" (before
#    "esponsible...") the same way.
$tailRange = $d.Range($afterPhrase, $afterRIdx)
$tailRange.Font.Size = $origSize + 1
$tailRange2 = $d.Range($afterPhrase, $afterRIdx)
$tailRange2.Font.Size = $origSize

Write-Output "Done"
